$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update 최종점수 (K) and MACRO_SCORE (N) columns for rows 2-7
$ws.Range("K2").Value = 56.5
$ws.Range("N2").Value = 53.62998959737769

$ws.Range("K3").Value = 56.5
$ws.Range("N3").Value = 53.62998959737769

$ws.Range("K4").Value = 52.3
$ws.Range("N4").Value = 53.62998959737769

$ws.Range("K5").Value = 47.3
$ws.Range("N5").Value = 53.62998959737769

$ws.Range("K6").Value = 45.1
$ws.Range("N6").Value = 53.62998959737769

$ws.Range("K7").Value = 44.5
$ws.Range("N7").Value = 53.62998959737769
